$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the title.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete() | Out-Null

# 2. Insert a new bold paragraph "Play Easter Island Free: A Fun, Engaging Slot
#    Game" right before the final paragraph (the one that used to hold the
#    image-generation prompt).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.InsertXML(
  "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
  "<w:r/>" +
  "<w:r><w:rPr><w:b/></w:rPr><w:t>Play Easter Island Free: A Fun, Engaging Slot Game</w:t></w:r>" +
  "</w:p>" +
  "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>"
) | Out-Null

# The InsertXML call above splits the former last paragraph in two, leaving a
# stray empty paragraph between the new text and the original last paragraph;
# remove that leftover paragraph mark.
$strayPara = $d.Paragraphs($count + 1)
$strayPara.Range.Delete() | Out-Null

# 3. Replace the old image-generation prompt text (now in the final paragraph)
#    with the new meta-description text, keeping its italic formatting intact.
$oldText = "Create a cartoon-style feature image for Easter Island that showcases a happy Maya warrior with glasses. The image should have a fun and vibrant feel, with the warrior holding a bunny or Easter egg to tie in with the Easter theme of the game. Use a colorful and engaging background, such as palm trees on a beach or a field of flowers, to make the image stand out. Make sure the Maya warrior has a big smile on their face to convey the fun and exciting atmosphere of the game."
$newText = "Experience the spring atmosphere of Easter Island with engaging gameplay and satisfying bonus features. Play for free now."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

Write-Output "Done: applied Easter Island meta-description relocation edit."
